$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing data right
$ws.Columns("D:D").Insert()

# Copy number formats/styles from (old D, now at) E into the new D column
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 37 and 79 are label-only rows with no data columns; remove the stray
# formatted cell created by the paste-format operation above.
$ws.Range("D37").Clear()
$ws.Range("D79").Clear()

# Populate the full data grid (columns D:L) with the restated financials
# Row 7
$ws.Range("D7").Value = 43281
$ws.Range("E7").Value = 42916
$ws.Range("F7").Value = 42551
$ws.Range("G7").Value = 42185
$ws.Range("H7").Value = 41820
$ws.Range("I7").Value = 41455
$ws.Range("J7").Value = 41090
$ws.Range("K7").Value = 40724
$ws.Range("L7").ClearContents()

# Row 8
$ws.Range("D8").Value = 170700
$ws.Range("E8").Value = 160400
$ws.Range("F8").Value = 166800
$ws.Range("G8").Value = 144300
$ws.Range("H8").Value = 124000
$ws.Range("I8").Value = 142300
$ws.Range("J8").Value = 120900
$ws.Range("K8").Value = 97700
$ws.Range("L8").ClearContents()

# Row 9
$ws.Range("D9").Value = 160900
$ws.Range("E9").Value = 156600
$ws.Range("F9").Value = 153300
$ws.Range("G9").Value = 132900
$ws.Range("H9").Value = 115200
$ws.Range("I9").Value = 112300
$ws.Range("J9").Value = 96000
$ws.Range("K9").Value = 83500
$ws.Range("L9").ClearContents()

# Row 10
$ws.Range("D10").Value = 9800
$ws.Range("E10").Value = 3800
$ws.Range("F10").Value = 13500
$ws.Range("G10").Value = 11400
$ws.Range("H10").Value = 8800
$ws.Range("I10").Value = 30000
$ws.Range("J10").Value = 24900
$ws.Range("K10").Value = 14200
$ws.Range("L10").ClearContents()

# Row 11
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("H11").ClearContents()
$ws.Range("I11").ClearContents()
$ws.Range("J11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()

# Row 12
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").ClearContents()

# Row 13
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").ClearContents()

# Row 14
$ws.Range("D14").Value = 600
$ws.Range("E14").Value = 1600
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 500
$ws.Range("H14").Value = 4200
$ws.Range("I14").Value = 16400
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 100
$ws.Range("L14").ClearContents()

# Row 15
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = "NA"
$ws.Range("L15").ClearContents()

# Row 16
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("G16").ClearContents()
$ws.Range("H16").ClearContents()
$ws.Range("I16").ClearContents()
$ws.Range("J16").ClearContents()
$ws.Range("K16").ClearContents()
$ws.Range("L16").ClearContents()

# Row 17
$ws.Range("D17").Value = 167100
$ws.Range("E17").Value = 162100
$ws.Range("F17").Value = 158600
$ws.Range("G17").Value = 137400
$ws.Range("H17").Value = 124700
$ws.Range("I17").Value = 134000
$ws.Range("J17").Value = 103600
$ws.Range("K17").Value = 89700
$ws.Range("L17").ClearContents()

# Row 18
$ws.Range("D18").Value = 3600
$ws.Range("E18").Value = -1700
$ws.Range("F18").Value = 8200
$ws.Range("G18").Value = 6900
$ws.Range("H18").Value = -700
$ws.Range("I18").Value = 8300
$ws.Range("J18").Value = 17300
$ws.Range("K18").Value = 8000
$ws.Range("L18").ClearContents()

# Row 19
$ws.Range("D19").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("G19").ClearContents()
$ws.Range("H19").ClearContents()
$ws.Range("I19").ClearContents()
$ws.Range("J19").ClearContents()
$ws.Range("K19").ClearContents()
$ws.Range("L19").ClearContents()

# Row 20
$ws.Range("D20").Value = -1100
$ws.Range("E20").Value = -600
$ws.Range("F20").Value = -400
$ws.Range("G20").Value = 500
$ws.Range("H20").Value = -900
$ws.Range("I20").Value = 2200
$ws.Range("J20").Value = 600
$ws.Range("K20").Value = 1600
$ws.Range("L20").ClearContents()

# Row 21
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "NA"
$ws.Range("G21").Value = "NA"
$ws.Range("H21").Value = "NA"
$ws.Range("I21").Value = 20400
$ws.Range("J21").Value = 26200
$ws.Range("K21").Value = "NA"
$ws.Range("L21").ClearContents()

# Row 22
$ws.Range("D22").Value = 300
$ws.Range("E22").Value = 300
$ws.Range("F22").Value = 300
$ws.Range("G22").Value = 700
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = "NA"
$ws.Range("K22").Value = 1600
$ws.Range("L22").ClearContents()

# Row 23
$ws.Range("D23").Value = 2200
$ws.Range("E23").Value = -2500
$ws.Range("F23").Value = 7500
$ws.Range("G23").Value = 6600
$ws.Range("H23").Value = -2500
$ws.Range("I23").Value = 9900
$ws.Range("J23").Value = 17900
$ws.Range("K23").Value = 8000
$ws.Range("L23").ClearContents()

# Row 24
$ws.Range("D24").Value = 1800
$ws.Range("E24").Value = -3500
$ws.Range("F24").Value = 3200
$ws.Range("G24").Value = 2000
$ws.Range("H24").Value = 1200
$ws.Range("I24").Value = 3100
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 2100
$ws.Range("L24").ClearContents()

# Row 25
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").ClearContents()

# Row 26
$ws.Range("D26").Value = 400
$ws.Range("E26").Value = 900
$ws.Range("F26").Value = 4200
$ws.Range("G26").Value = 4700
$ws.Range("H26").Value = -3800
$ws.Range("I26").Value = 6900
$ws.Range("J26").Value = 17300
$ws.Range("K26").Value = 5900
$ws.Range("L26").ClearContents()

# Row 27
$ws.Range("D27").Value = 400
$ws.Range("E27").Value = 900
$ws.Range("F27").Value = 4200
$ws.Range("G27").Value = 4600
$ws.Range("H27").Value = -3100
$ws.Range("I27").Value = 4100
$ws.Range("J27").Value = 14900
$ws.Range("K27").Value = 4700
$ws.Range("L27").ClearContents()

# Row 28
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").ClearContents()

# Row 29
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = -3000
$ws.Range("K29").Value = -25100
$ws.Range("L29").ClearContents()

# Row 30
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").ClearContents()

# Row 31
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").ClearContents()

# Row 32
$ws.Range("D32").Value = 1100
$ws.Range("E32").Value = 600
$ws.Range("F32").Value = 400
$ws.Range("G32").Value = -500
$ws.Range("H32").Value = 900
$ws.Range("I32").Value = -2200
$ws.Range("J32").Value = -600
$ws.Range("K32").Value = -1600
$ws.Range("L32").ClearContents()

# Row 33
$ws.Range("D33").Value = 400
$ws.Range("E33").Value = 900
$ws.Range("F33").Value = 4200
$ws.Range("G33").Value = 4600
$ws.Range("H33").Value = -3100
$ws.Range("I33").Value = 4100
$ws.Range("J33").Value = 11900
$ws.Range("K33").Value = -20400
$ws.Range("L33").ClearContents()

# Row 34
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()

# Row 35
$ws.Range("D35").Value = 400
$ws.Range("E35").Value = 900
$ws.Range("F35").Value = 4200
$ws.Range("G35").Value = 4600
$ws.Range("H35").Value = -3100
$ws.Range("I35").Value = 4100
$ws.Range("J35").Value = 11900
$ws.Range("K35").Value = -20400
$ws.Range("L35").ClearContents()

# Row 38
$ws.Range("D38").Value = 43281
$ws.Range("E38").Value = 42916
$ws.Range("F38").Value = 42551
$ws.Range("G38").Value = 42185
$ws.Range("H38").Value = 41820
$ws.Range("I38").Value = 41455
$ws.Range("J38").Value = 41090
$ws.Range("K38").Value = 40724
$ws.Range("L38").ClearContents()

# Row 39
$ws.Range("D39").ClearContents()
$ws.Range("E39").ClearContents()
$ws.Range("F39").ClearContents()
$ws.Range("G39").ClearContents()
$ws.Range("H39").ClearContents()
$ws.Range("I39").ClearContents()
$ws.Range("J39").ClearContents()
$ws.Range("K39").ClearContents()
$ws.Range("L39").ClearContents()

# Row 40
$ws.Range("D40").ClearContents()
$ws.Range("E40").ClearContents()
$ws.Range("F40").ClearContents()
$ws.Range("G40").ClearContents()
$ws.Range("H40").ClearContents()
$ws.Range("I40").ClearContents()
$ws.Range("J40").ClearContents()
$ws.Range("K40").ClearContents()
$ws.Range("L40").ClearContents()

# Row 41
$ws.Range("D41").Value = 20700
$ws.Range("E41").Value = 17400
$ws.Range("F41").Value = 24100
$ws.Range("G41").Value = 22200
$ws.Range("H41").Value = 14300
$ws.Range("I41").Value = 25800
$ws.Range("J41").Value = 20500
$ws.Range("K41").Value = 18400
$ws.Range("L41").ClearContents()

# Row 42
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").ClearContents()

# Row 43
$ws.Range("D43").Value = 5400
$ws.Range("E43").Value = 7800
$ws.Range("F43").Value = 4500
$ws.Range("G43").Value = 7100
$ws.Range("H43").Value = 7000
$ws.Range("I43").Value = 11700
$ws.Range("J43").Value = 4200
$ws.Range("K43").Value = 8300
$ws.Range("L43").ClearContents()

# Row 44
$ws.Range("D44").Value = 16000
$ws.Range("E44").Value = 12400
$ws.Range("F44").Value = 11000
$ws.Range("G44").Value = 11600
$ws.Range("H44").Value = 10100
$ws.Range("I44").Value = 19000
$ws.Range("J44").Value = 7300
$ws.Range("K44").Value = 8700
$ws.Range("L44").ClearContents()

# Row 45
$ws.Range("D45").Value = 800
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 1600
$ws.Range("G45").Value = 900
$ws.Range("H45").Value = 800
$ws.Range("I45").Value = 400
$ws.Range("J45").Value = 400
$ws.Range("K45").Value = 800
$ws.Range("L45").ClearContents()

# Row 46
$ws.Range("D46").Value = 42900
$ws.Range("E46").Value = 37600
$ws.Range("F46").Value = 41200
$ws.Range("G46").Value = 41700
$ws.Range("H46").Value = 32200
$ws.Range("I46").Value = 41400
$ws.Range("J46").Value = 32300
$ws.Range("K46").Value = 36100
$ws.Range("L46").ClearContents()

# Row 47
$ws.Range("D47").Value = 18700
$ws.Range("E47").Value = 16200
$ws.Range("F47").Value = 7000
$ws.Range("G47").Value = 6400
$ws.Range("H47").Value = 8200
$ws.Range("I47").Value = 21100
$ws.Range("J47").Value = 23400
$ws.Range("K47").Value = 11300
$ws.Range("L47").ClearContents()

# Row 48
$ws.Range("D48").Value = 99600
$ws.Range("E48").Value = 102600
$ws.Range("F48").Value = 109700
$ws.Range("G48").Value = 116400
$ws.Range("H48").Value = 120300
$ws.Range("I48").Value = 240800
$ws.Range("J48").Value = 112500
$ws.Range("K48").Value = 133300
$ws.Range("L48").ClearContents()

# Row 49
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").ClearContents()

# Row 50
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").ClearContents()

# Row 51
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").ClearContents()

# Row 52
$ws.Range("D52").Value = 600
$ws.Range("E52").Value = 300
$ws.Range("F52").Value = 7900
$ws.Range("G52").Value = 7000
$ws.Range("H52").Value = 6500
$ws.Range("I52").Value = 200
$ws.Range("J52").Value = 2600
$ws.Range("K52").Value = 5600
$ws.Range("L52").ClearContents()

# Row 53
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").ClearContents()

# Row 54
$ws.Range("D54").Value = 161800
$ws.Range("E54").Value = 156800
$ws.Range("F54").Value = 165800
$ws.Range("G54").Value = 171600
$ws.Range("H54").Value = 167300
$ws.Range("I54").Value = 183100
$ws.Range("J54").Value = 170800
$ws.Range("K54").Value = 162100
$ws.Range("L54").ClearContents()

# Row 55
$ws.Range("D55").ClearContents()
$ws.Range("E55").ClearContents()
$ws.Range("F55").ClearContents()
$ws.Range("G55").ClearContents()
$ws.Range("H55").ClearContents()
$ws.Range("I55").ClearContents()
$ws.Range("J55").ClearContents()
$ws.Range("K55").ClearContents()
$ws.Range("L55").ClearContents()

# Row 56
$ws.Range("D56").ClearContents()
$ws.Range("E56").ClearContents()
$ws.Range("F56").ClearContents()
$ws.Range("G56").ClearContents()
$ws.Range("H56").ClearContents()
$ws.Range("I56").ClearContents()
$ws.Range("J56").ClearContents()
$ws.Range("K56").ClearContents()
$ws.Range("L56").ClearContents()

# Row 57
$ws.Range("D57").Value = 15600
$ws.Range("E57").Value = 13800
$ws.Range("F57").Value = 19800
$ws.Range("G57").Value = 17700
$ws.Range("H57").Value = 14500
$ws.Range("I57").Value = 15200
$ws.Range("J57").Value = 15600
$ws.Range("K57").Value = 23400
$ws.Range("L57").ClearContents()

# Row 58
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("F58").Value = 200
$ws.Range("G58").Value = 1700
$ws.Range("H58").Value = 5000
$ws.Range("I58").Value = 3300
$ws.Range("J58").Value = 2100
$ws.Range("K58").Value = 5600
$ws.Range("L58").ClearContents()

# Row 59
$ws.Range("D59").Value = 6400
$ws.Range("E59").Value = 3800
$ws.Range("F59").Value = 1500
$ws.Range("G59").Value = 1400
$ws.Range("H59").Value = 1500
$ws.Range("I59").Value = 100
$ws.Range("J59").Value = 200
$ws.Range("K59").Value = 0
$ws.Range("L59").ClearContents()

# Row 60
$ws.Range("D60").Value = 22000
$ws.Range("E60").Value = 17600
$ws.Range("F60").Value = 21500
$ws.Range("G60").Value = 20800
$ws.Range("H60").Value = 21100
$ws.Range("I60").Value = 16800
$ws.Range("J60").Value = 17900
$ws.Range("K60").Value = 29000
$ws.Range("L60").ClearContents()

# Row 61
$ws.Range("D61").Value = 1000
$ws.Range("E61").Value = 1200
$ws.Range("F61").Value = 1200
$ws.Range("G61").Value = 1300
$ws.Range("H61").Value = 5200
$ws.Range("I61").Value = 9800
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2900
$ws.Range("L61").ClearContents()

# Row 62
$ws.Range("D62").Value = 51900
$ws.Range("E62").Value = 48700
$ws.Range("F62").Value = 51300
$ws.Range("G62").Value = 44600
$ws.Range("H62").Value = 39500
$ws.Range("I62").Value = 43400
$ws.Range("J62").Value = 40900
$ws.Range("K62").Value = 43800
$ws.Range("L62").ClearContents()

# Row 63
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").ClearContents()

# Row 64
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").ClearContents()

# Row 65
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").ClearContents()

# Row 66
$ws.Range("D66").Value = 74900
$ws.Range("E66").Value = 67500
$ws.Range("F66").Value = 74000
$ws.Range("G66").Value = 66700
$ws.Range("H66").Value = 81700
$ws.Range("I66").Value = 85300
$ws.Range("J66").Value = 68200
$ws.Range("K66").Value = 73800
$ws.Range("L66").ClearContents()

# Row 67
$ws.Range("D67").ClearContents()
$ws.Range("E67").ClearContents()
$ws.Range("F67").ClearContents()
$ws.Range("G67").ClearContents()
$ws.Range("H67").ClearContents()
$ws.Range("I67").ClearContents()
$ws.Range("J67").ClearContents()
$ws.Range("K67").ClearContents()
$ws.Range("L67").ClearContents()

# Row 68
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()

# Row 69
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").ClearContents()

# Row 70
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").ClearContents()

# Row 71
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()

# Row 72
$ws.Range("D72").Value = -199500
$ws.Range("E72").Value = -197100
$ws.Range("F72").Value = -204100
$ws.Range("G72").Value = -191000
$ws.Range("H72").Value = -206500
$ws.Range("I72").Value = -197900
$ws.Range("J72").Value = -191700
$ws.Range("K72").Value = -218800
$ws.Range("L72").ClearContents()

# Row 73
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").ClearContents()

# Row 74
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()

# Row 75
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").ClearContents()

# Row 76
$ws.Range("D76").Value = 86900
$ws.Range("E76").Value = 89300
$ws.Range("F76").Value = 91800
$ws.Range("G76").Value = 104900
$ws.Range("H76").Value = 85600
$ws.Range("I76").Value = 97800
$ws.Range("J76").Value = 102600
$ws.Range("K76").Value = 88300
$ws.Range("L76").ClearContents()

# Row 77
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()

# Row 80
$ws.Range("D80").Value = 43281
$ws.Range("E80").Value = 42916
$ws.Range("F80").Value = 42551
$ws.Range("G80").Value = 42185
$ws.Range("H80").Value = 41820
$ws.Range("I80").Value = 41455
$ws.Range("J80").Value = 41090
$ws.Range("K80").Value = 40724
$ws.Range("L80").ClearContents()

# Row 81
$ws.Range("D81").Value = 400
$ws.Range("E81").Value = 900
$ws.Range("F81").Value = 4200
$ws.Range("G81").Value = 4600
$ws.Range("H81").Value = -3100
$ws.Range("I81").Value = 4100
$ws.Range("J81").Value = 11900
$ws.Range("K81").Value = -20400
$ws.Range("L81").ClearContents()

# Row 82
$ws.Range("D82").ClearContents()
$ws.Range("E82").ClearContents()
$ws.Range("F82").ClearContents()
$ws.Range("G82").ClearContents()
$ws.Range("H82").ClearContents()
$ws.Range("I82").ClearContents()
$ws.Range("J82").ClearContents()
$ws.Range("K82").ClearContents()
$ws.Range("L82").ClearContents()

# Row 83
$ws.Range("D83").Value = "NA"
$ws.Range("E83").Value = "NA"
$ws.Range("F83").Value = "NA"
$ws.Range("G83").Value = "NA"
$ws.Range("H83").Value = "NA"
$ws.Range("I83").Value = 9900
$ws.Range("J83").Value = 8300
$ws.Range("K83").Value = "NA"
$ws.Range("L83").ClearContents()

# Row 84
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").ClearContents()

# Row 85
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").ClearContents()

# Row 86
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()

# Row 87
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").ClearContents()

# Row 88
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()

# Row 89
$ws.Range("D89").Value = 16000
$ws.Range("E89").Value = 3500
$ws.Range("F89").Value = 28500
$ws.Range("G89").Value = 19400
$ws.Range("H89").Value = 5500
$ws.Range("I89").Value = 34400
$ws.Range("J89").Value = 42600
$ws.Range("K89").Value = 22900
$ws.Range("L89").ClearContents()

# Row 90
$ws.Range("D90").ClearContents()
$ws.Range("E90").ClearContents()
$ws.Range("F90").ClearContents()
$ws.Range("G90").ClearContents()
$ws.Range("H90").ClearContents()
$ws.Range("I90").ClearContents()
$ws.Range("J90").ClearContents()
$ws.Range("K90").ClearContents()
$ws.Range("L90").ClearContents()

# Row 91
$ws.Range("D91").Value = -8600
$ws.Range("E91").Value = -7600
$ws.Range("F91").Value = -6800
$ws.Range("G91").Value = -6200
$ws.Range("H91").Value = -10900
$ws.Range("I91").Value = -26200
$ws.Range("J91").Value = -22800
$ws.Range("K91").Value = -22500
$ws.Range("L91").ClearContents()

# Row 92
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").ClearContents()

# Row 93
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").ClearContents()

# Row 94
$ws.Range("D94").Value = -9600
$ws.Range("E94").Value = -6600
$ws.Range("F94").Value = -7300
$ws.Range("G94").Value = -2600
$ws.Range("H94").Value = -11800
$ws.Range("I94").Value = -29400
$ws.Range("J94").Value = -28300
$ws.Range("K94").Value = "NA"
$ws.Range("L94").ClearContents()

# Row 95
$ws.Range("D95").ClearContents()
$ws.Range("E95").ClearContents()
$ws.Range("F95").ClearContents()
$ws.Range("G95").ClearContents()
$ws.Range("H95").ClearContents()
$ws.Range("I95").ClearContents()
$ws.Range("J95").ClearContents()
$ws.Range("K95").ClearContents()
$ws.Range("L95").ClearContents()

# Row 96
$ws.Range("D96").Value = -2900
$ws.Range("E96").Value = -3500
$ws.Range("F96").Value = -17300
$ws.Range("G96").Value = -500
$ws.Range("H96").Value = -3600
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = -1400
$ws.Range("L96").ClearContents()

# Row 97
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()

# Row 98
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").ClearContents()

# Row 99
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()

# Row 100
$ws.Range("D100").Value = -3100
$ws.Range("E100").Value = -3600
$ws.Range("F100").Value = -19300
$ws.Range("G100").Value = -8900
$ws.Range("H100").Value = -5300
$ws.Range("I100").Value = 400
$ws.Range("J100").Value = -11600
$ws.Range("K100").Value = "NA"
$ws.Range("L100").ClearContents()

# Row 101
$ws.Range("D101").Value = "NA"
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = "NA"
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = "NA"
$ws.Range("L101").ClearContents()

# Row 102
$ws.Range("D102").Value = 3300
$ws.Range("E102").Value = -6700
$ws.Range("F102").Value = 1900
$ws.Range("G102").Value = 7900
$ws.Range("H102").Value = -11500
$ws.Range("I102").Value = 5400
$ws.Range("J102").Value = 2700
$ws.Range("K102").Value = 5000
$ws.Range("L102").ClearContents()
